$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column I = DAMSLTag, Column J = DialogAct
$ws.Cells.Item(17, 9).Value = "sd"
$ws.Cells.Item(17, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(25, 9).Value = "sd"
$ws.Cells.Item(25, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(27, 9).Value = "sd"
$ws.Cells.Item(27, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(28, 9).Value = "sv"
$ws.Cells.Item(28, 10).Value = "Statement-opinion"
$ws.Cells.Item(30, 9).Value = "ba"
$ws.Cells.Item(30, 10).Value = "Appreciation"
$ws.Cells.Item(36, 9).Value = "ba"
$ws.Cells.Item(36, 10).Value = "Appreciation"
$ws.Cells.Item(39, 9).Value = "sd"
$ws.Cells.Item(39, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(40, 9).Value = "b"
$ws.Cells.Item(40, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(41, 9).Value = "sv"
$ws.Cells.Item(41, 10).Value = "Statement-opinion"
$ws.Cells.Item(56, 9).Value = "sv"
$ws.Cells.Item(56, 10).Value = "Statement-opinion"
$ws.Cells.Item(57, 9).Value = "aa"
$ws.Cells.Item(57, 10).Value = "Agree/Accept"
$ws.Cells.Item(59, 9).Value = "%"
$ws.Cells.Item(59, 10).Value = "Uninterpretable"
$ws.Cells.Item(60, 9).Value = "%"
$ws.Cells.Item(60, 10).Value = "Uninterpretable"
$ws.Cells.Item(87, 9).Value = "aa"
$ws.Cells.Item(87, 10).Value = "Agree/Accept"
$ws.Cells.Item(111, 9).Value = "sd"
$ws.Cells.Item(111, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(112, 9).Value = "%"
$ws.Cells.Item(112, 10).Value = "Uninterpretable"
$ws.Cells.Item(117, 9).Value = "aa"
$ws.Cells.Item(117, 10).Value = "Agree/Accept"
$ws.Cells.Item(118, 9).Value = "%"
$ws.Cells.Item(118, 10).Value = "Uninterpretable"
$ws.Cells.Item(121, 9).Value = "sd"
$ws.Cells.Item(121, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(137, 9).Value = "sd"
$ws.Cells.Item(137, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(140, 9).Value = "ba"
$ws.Cells.Item(140, 10).Value = "Appreciation"
$ws.Cells.Item(142, 9).Value = "sd"
$ws.Cells.Item(142, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(143, 9).Value = "sd"
$ws.Cells.Item(143, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(144, 9).Value = "sd"
$ws.Cells.Item(144, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(149, 9).Value = "sd"
$ws.Cells.Item(149, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(151, 9).Value = "b"
$ws.Cells.Item(151, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(153, 9).Value = "b"
$ws.Cells.Item(153, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(158, 9).Value = "%"
$ws.Cells.Item(158, 10).Value = "Uninterpretable"
$ws.Cells.Item(174, 9).Value = "sd"
$ws.Cells.Item(174, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(180, 9).Value = "aa"
$ws.Cells.Item(180, 10).Value = "Agree/Accept"
$ws.Cells.Item(196, 9).Value = "sd"
$ws.Cells.Item(196, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(215, 9).Value = "sd"
$ws.Cells.Item(215, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(217, 9).Value = "aa"
$ws.Cells.Item(217, 10).Value = "Agree/Accept"
$ws.Cells.Item(219, 9).Value = "aa"
$ws.Cells.Item(219, 10).Value = "Agree/Accept"
$ws.Cells.Item(225, 9).Value = "sd"
$ws.Cells.Item(225, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(234, 9).Value = "sd"
$ws.Cells.Item(234, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(236, 9).Value = "aa"
$ws.Cells.Item(236, 10).Value = "Agree/Accept"
$ws.Cells.Item(238, 9).Value = "sd"
$ws.Cells.Item(238, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(247, 9).Value = "sv"
$ws.Cells.Item(247, 10).Value = "Statement-opinion"
$ws.Cells.Item(266, 9).Value = "aa"
$ws.Cells.Item(266, 10).Value = "Agree/Accept"
$ws.Cells.Item(267, 9).Value = "aa"
$ws.Cells.Item(267, 10).Value = "Agree/Accept"
$ws.Cells.Item(268, 9).Value = "sd"
$ws.Cells.Item(268, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(302, 9).Value = "sd"
$ws.Cells.Item(302, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(308, 9).Value = "ba"
$ws.Cells.Item(308, 10).Value = "Appreciation"
$ws.Cells.Item(309, 9).Value = "sv"
$ws.Cells.Item(309, 10).Value = "Statement-opinion"
$ws.Cells.Item(315, 9).Value = "sd"
$ws.Cells.Item(315, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(320, 9).Value = "b"
$ws.Cells.Item(320, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(321, 9).Value = "sd"
$ws.Cells.Item(321, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(334, 9).Value = "sv"
$ws.Cells.Item(334, 10).Value = "Statement-opinion"
$ws.Cells.Item(342, 9).Value = "sv"
$ws.Cells.Item(342, 10).Value = "Statement-opinion"
$ws.Cells.Item(347, 9).Value = "sd"
$ws.Cells.Item(347, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(356, 9).Value = "sd"
$ws.Cells.Item(356, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(358, 9).Value = "sd"
$ws.Cells.Item(358, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(364, 9).Value = "b"
$ws.Cells.Item(364, 10).Value = "Acknowledge (Backchannel)"
